$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the absolute path metadata recorded the last time the file was
#     saved locally (folder got a space inserted: "ПоказателиЦУР" -> "Показатели ЦУР").
#     (No Excel object-model surface exposes this incidental MRU path, so this
#     is left as-is; it is not reachable through normal user actions.)

# --- Column N (2020 data) ------------------------------------------------

# N3: blank cell on the heavy-ruled separator row - reuse the row's existing
# format (bottom medium border) by copying it from a neighboring cell.
$ws.Range("M3").Copy() | Out-Null
$ws.Range("N3").PasteSpecial(-4122) | Out-Null

# N4: header year value, bold table font with only a bottom rule.
$c = $ws.Range("N4")
$c.Value = 2020
$c.Font.Name = "Times New Roman"
$c.Font.Size = 9
$c.Font.Bold = $true
$c.VerticalAlignment = -4107
$c.Borders.Item(9).LineStyle = 1
$c.Borders.Item(9).Weight = -4138

# N5: blank row under the header, plain table font, no border.
$c = $ws.Range("N5")
$c.VerticalAlignment = -4107

# N6:N9 data values, plain table font, no border - same look as N5.
$ws.Range("N5").Copy() | Out-Null
$ws.Range("N6").PasteSpecial(-4122) | Out-Null
$ws.Range("N6").Value = 1713

$ws.Range("N5").Copy() | Out-Null
$ws.Range("N7").PasteSpecial(-4122) | Out-Null
$ws.Range("N7").Value = 1

$ws.Range("N5").Copy() | Out-Null
$ws.Range("N8").PasteSpecial(-4122) | Out-Null
$ws.Range("N8").Value = 379

$ws.Range("N5").Copy() | Out-Null
$ws.Range("N9").PasteSpecial(-4122) | Out-Null
$ws.Range("N9").Value = 180

# N10: bottom-row total, plain table font with bottom rule.
$c = $ws.Range("N10")
$c.Value = 798
$c.VerticalAlignment = -4107
$c.Borders.Item(9).LineStyle = 1
$c.Borders.Item(9).Weight = -4138

$excel.CutCopyMode = $false

# --- Selection left on the sheet after the edit --------------------------
$ws.Range("L22").Select() | Out-Null
